$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 54
$ws.Range("R2").Value = 1500
$ws.Range("S2").Value = 1620

# Row 3 updates
$ws.Range("B3").Value = 56
$ws.Range("D3").Value = 28
$ws.Range("G3").Value = 28
$ws.Range("R3").Value = 1680
$ws.Range("T3").Value = 840
$ws.Range("V3").Value = 840
